# The commit permutes the data rows 6-13 on the active sheet: each row's
# entire content (columns A..AY) moves to a different row number, while the
# row numbers themselves (and everything outside rows 6-13) stay put.
#
# Mapping of destination row -> source row (i.e. destination row R ends up
# holding what source row mapping[R] held before the edit):
#   6 <- 7
#   7 <- 9
#   8 <- 12
#   9 <- 13
#   10 <- 11
#   11 <- 10
#   12 <- 6
#   13 <- 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 6
$lastRow = 13
$firstCol = 1
$lastCol = 51   # column AY, matches the sheet's used range / dimension

$mapping = @{6=7; 7=9; 8=12; 9=13; 10=11; 11=10; 12=6; 13=8}

# 1) Snapshot every cell in rows 6-13 (all used columns) before touching
#    anything, so reads never see values we've already overwritten.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# 2) Write the permuted values back out.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $value = $snapshot["$src,$c"]
        $cell = $ws.Cells.Item($r, $c)

        if ($value -ne $null -and $value.GetType().Name -eq "String") {
            # Guard against Excel auto-converting date-looking strings
            # (e.g. "2023-09-03") into real date serials on assignment.
            $cell.NumberFormat = "@"
        }

        $cell.Value2 = $value
    }
}
